$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-three character used in row 40 price (U+2083)
$sub3 = [char]0x2083

# Force the Price column (D) cells being updated to stay text, so values like
# "589.30", "1.00" or "50.00" keep their exact digits/trailing zeros instead of
# being auto-converted to numbers by Excel.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '71.803.31'
$ws.Cells.Item(2, 5).Value = '  +3.64%  '
$ws.Cells.Item(3, 4).Value = '3.694.91'
$ws.Cells.Item(3, 5).Value = '  +9.09%  '
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(5, 4).Value = '589.30'
$ws.Cells.Item(5, 5).Value = '  +1.53%  '
$ws.Cells.Item(6, 4).Value = '179.12'
$ws.Cells.Item(6, 5).Value = '  +0.24%  '
$ws.Cells.Item(7, 4).Value = '3.688.62'
$ws.Cells.Item(7, 5).Value = '  +9.11%  '
$ws.Cells.Item(8, 5).Value = '  +5.44%  '
$ws.Cells.Item(9, 5).Value = '  +0.05%  '
$ws.Cells.Item(10, 4).Value = '0.202'
$ws.Cells.Item(10, 5).Value = '  +1.88%  '
$ws.Cells.Item(11, 5).Value = '  +4.88%  '
$ws.Cells.Item(12, 4).Value = '50.00'
$ws.Cells.Item(12, 5).Value = '  +3.45%  '
$ws.Cells.Item(13, 5).Value = '  +1.68%  '
$ws.Cells.Item(14, 4).Value = '4.279.98'
$ws.Cells.Item(14, 5).Value = '  +8.84%  '
$ws.Cells.Item(15, 4).Value = '681.66'
$ws.Cells.Item(15, 5).Value = '  -0.17%  '
$ws.Cells.Item(16, 4).Value = '8.99'
$ws.Cells.Item(16, 5).Value = '  +4.90%  '
$ws.Cells.Item(17, 4).Value = '71.926.40'
$ws.Cells.Item(17, 5).Value = '  +3.66%  '
$ws.Cells.Item(18, 4).Value = '3.665.59'
$ws.Cells.Item(18, 5).Value = '  +8.23%  '
$ws.Cells.Item(19, 5).Value = '  +2.03%  '
$ws.Cells.Item(20, 4).Value = '18.06'
$ws.Cells.Item(20, 5).Value = '  +2.40%  '
$ws.Cells.Item(21, 4).Value = '11.66'
$ws.Cells.Item(21, 5).Value = '  +3.46%  '
$ws.Cells.Item(22, 4).Value = '0.940'
$ws.Cells.Item(22, 5).Value = '  +3.49%  '
$ws.Cells.Item(23, 4).Value = '6.22'
$ws.Cells.Item(23, 5).Value = '  +16.03%  '
$ws.Cells.Item(24, 4).Value = '17.82'
$ws.Cells.Item(24, 5).Value = '  +4.11%  '
$ws.Cells.Item(25, 4).Value = '103.48'
$ws.Cells.Item(25, 5).Value = '  +2.31%  '
$ws.Cells.Item(26, 5).Value = '  +3.91%  '
$ws.Cells.Item(27, 4).Value = '2.85'
$ws.Cells.Item(27, 5).Value = '  +5.56%  '
$ws.Cells.Item(28, 4).Value = '10.18'
$ws.Cells.Item(28, 5).Value = '  +5.01%  '
$ws.Cells.Item(29, 4).Value = '35.52'
$ws.Cells.Item(29, 5).Value = '  +6.17%  '
$ws.Cells.Item(30, 4).Value = '9.20'
$ws.Cells.Item(30, 5).Value = '  +5.40%  '
$ws.Cells.Item(31, 5).Value = '  +6.23%  '
$ws.Cells.Item(32, 5).Value = '  +9.45%  '
$ws.Cells.Item(33, 4).Value = '575.58'
$ws.Cells.Item(33, 5).Value = '  +4.77%  '
$ws.Cells.Item(34, 4).Value = '11.31'
$ws.Cells.Item(34, 5).Value = '  +2.54%  '
$ws.Cells.Item(35, 5).Value = '  +3.96%  '
$ws.Cells.Item(36, 4).Value = '59.61'
$ws.Cells.Item(36, 5).Value = '  +3.12%  '
$ws.Cells.Item(37, 4).Value = '3.765.74'
$ws.Cells.Item(37, 5).Value = '  +4.58%  '
$ws.Cells.Item(38, 4).Value = '1.00'
$ws.Cells.Item(39, 5).Value = '  +3.62%  '
$ws.Cells.Item(40, 4).Value = [string]::Concat('0.0', $sub3, '0778')
$ws.Cells.Item(40, 5).Value = '  +4.64%  '
$ws.Cells.Item(41, 4).Value = '35.46'
$ws.Cells.Item(41, 5).Value = '  +0.31%  '
$ws.Cells.Item(42, 5).Value = '  +5.25%  '
$ws.Cells.Item(43, 4).Value = '0.0464'
$ws.Cells.Item(43, 5).Value = '  +9.22%  '
$ws.Cells.Item(44, 5).Value = '  +3.13%  '
$ws.Cells.Item(45, 4).Value = '0.350'
$ws.Cells.Item(45, 5).Value = '  +4.58%  '
$ws.Cells.Item(46, 4).Value = '2.90'
$ws.Cells.Item(46, 5).Value = '  +8.65%  '
$ws.Cells.Item(47, 4).Value = '3.38'
$ws.Cells.Item(47, 5).Value = '  +0.08%  '
$ws.Cells.Item(48, 5).Value = '  +4.21%  '
$ws.Cells.Item(49, 5).Value = '  +2.60%  '
$ws.Cells.Item(51, 4).Value = '134.22'
$ws.Cells.Item(51, 5).Value = '  +3.25%  '
